$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row "005206566 / LEVI / 50091" right before the existing
#    row for 004376145 / LUCYENE (spreadsheet row 7).
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "'005206566"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = "LEVI"
$ws.Range("C7").Value = 50091

# 2) Remove the row for 002277249 / DANILO / 2239.1 (now at row 43 after
#    the insertion above shifted everything below row 7 down by one).
$ws.Rows.Item(43).Delete()

# 3) Remove the old row for 005206566 / LEVI / 91 (now at row 105 after
#    the insertion above shifted it down by one, and it is below the
#    DANILO row removed above so its own row number is unaffected by
#    that deletion).
$ws.Rows.Item(105).Delete()
